# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.715.62"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.338.57"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.86"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.68"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.32"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.84"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.753.37"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.685.35"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.352.04"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.43"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "325.84"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.24"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("E24").Value = "  +12.51%  "
$ws.Range("E25").Value = "  +3.96%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.31"
$ws.Range("E27").Value = "  +7.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.43"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0726"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.40"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.889"
$ws.Range("E37").Value = "  -5.37%  "
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.57"
$ws.Range("E39").Value = "  +3.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.93"
$ws.Range("E40").Value = "  +9.69%  "
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "281.48"
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.11"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0500"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.559"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.37"
$ws.Range("E48").Value = "  +7.87%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.11"
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("E51").Value = "  +1.27%  "
